$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are stored as text so values like "1.000" or
# "27.922.83" keep their exact formatting instead of being auto-converted
# to numbers by Excel.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "27.922.83"
$ws.Range("E2").Value = "  +0.57%  "

# Row 3
$ws.Range("D3").Value = "1.767.49"
$ws.Range("E3").Value = "  -0.48%  "

# Row 5
$ws.Range("D5").Value = "328.81"
$ws.Range("E5").Value = "  +0.48%  "

# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").Value = "0.4783"
$ws.Range("E7").Value = "  +4.34%  "

# Row 8
$ws.Range("D8").Value = "0.3532"
$ws.Range("E8").Value = "  -1.43%  "

# Row 9
$ws.Range("D9").Value = "43.35"
$ws.Range("E9").Value = "  +3.74%  "

# Row 10
$ws.Range("D10").Value = "0.07397"
$ws.Range("E10").Value = "  -1.29%  "

# Row 11
$ws.Range("E11").Value = "  -1.81%  "

# Row 12
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.02%  "

# Row 13
$ws.Range("D13").Value = "20.64"
$ws.Range("E13").Value = "  -0.91%  "

# Row 14
$ws.Range("D14").Value = "6.019"
$ws.Range("E14").Value = "  -0.34%  "

# Row 15
$ws.Range("D15").Value = "7.181"
$ws.Range("E15").Value = "  -0.37%  "

# Row 16
$ws.Range("D16").Value = "1.767.65"
$ws.Range("E16").Value = "  -0.66%  "

# Row 17
$ws.Range("D17").Value = "92.27"
$ws.Range("E17").Value = "  -1.50%  "

# Row 18
$ws.Range("E18").Value = "  -0.42%  "

# Row 19
$ws.Range("D19").Value = "0.06422"
$ws.Range("E19").Value = "  -0.12%  "

# Row 20
$ws.Range("E20").Value = "  -0.02%  "

# Row 21
$ws.Range("E21").Value = "  -0.84%  "

# Row 22
$ws.Range("D22").Value = "5.800"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").Value = "27.962.37"
$ws.Range("E23").Value = "  +0.59%  "

# Row 24
$ws.Range("E24").Value = "  -1.46%  "

# Row 25
$ws.Range("D25").Value = "2.157"
$ws.Range("E25").Value = "  +3.71%  "

# Row 26
$ws.Range("D26").Value = "164.90"
$ws.Range("E26").Value = "  +0.70%  "

# Row 27
$ws.Range("D27").Value = "20.09"
$ws.Range("E27").Value = "  -0.95%  "

# Row 28
$ws.Range("D28").Value = "1.971.51"
$ws.Range("E28").Value = "  -0.37%  "

# Row 29
$ws.Range("D29").Value = "2.214"
$ws.Range("E29").Value = "  +1.27%  "

# Row 30
$ws.Range("D30").Value = "123.49"
$ws.Range("E30").Value = "  -1.80%  "

# Row 31
$ws.Range("D31").Value = "1.075"
$ws.Range("E31").Value = "  -2.81%  "

# Row 32
$ws.Range("D32").Value = "0.09417"
$ws.Range("E32").Value = "  +2.12%  "

# Row 33
$ws.Range("D33").Value = "3.661"
$ws.Range("E33").Value = "  -0.18%  "

# Row 34
$ws.Range("D34").Value = "5.558"
$ws.Range("E34").Value = "  +0.35%  "

# Row 35
$ws.Range("E35").Value = "  -1.47%  "

# Row 36
$ws.Range("D36").Value = "0.06111"
$ws.Range("E36").Value = "  -0.84%  "

# Row 37
$ws.Range("D37").Value = "0.02261"
$ws.Range("E37").Value = "  -1.35%  "

# Row 38
$ws.Range("E38").Value = "  -1.03%  "

# Row 39
$ws.Range("D39").Value = "4.903"
$ws.Range("E39").Value = "  -1.15%  "

# Row 40
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "1.451"
$ws.Range("E40").Value = "  +4.57%  "

# Row 41
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6177"
$ws.Range("E41").Value = "  -2.37%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.189"
$ws.Range("E42").Value = "  +0.06%  "

# Row 43
$ws.Range("D43").Value = "7.753"
$ws.Range("E43").Value = "  -1.22%  "

# Row 44
$ws.Range("D44").Value = "13.08"
$ws.Range("E44").Value = "  -1.89%  "

# Row 45
$ws.Range("D45").Value = "3.746"
$ws.Range("E45").Value = "  +0.11%  "

# Row 46
$ws.Range("D46").Value = "0.5806"
$ws.Range("E46").Value = "  -2.03%  "

# Row 47
$ws.Range("D47").Value = "124.03"
$ws.Range("E47").Value = "  +1.00%  "

# Row 49
$ws.Range("D49").Value = "1.130"
$ws.Range("E49").Value = "  -0.64%  "

# Row 50
$ws.Range("D50").Value = "0.06807"
$ws.Range("E50").Value = "  -1.68%  "

# Row 51
$ws.Range("D51").Value = "72.10"
$ws.Range("E51").Value = "  -0.40%  "
